$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 458, shifting existing rows 458-477 down to 459-478.
$ws.Rows("458:458").Insert(4)

# Populate the newly inserted row 458 with the new record.
$ws.Range("A458").Value = 3
$ws.Range("B458").Value = "Femacal de La Calera"
$ws.Range("C458").Value = "Coquimbo"
$ws.Range("D458").Value = 44939
$ws.Range("E458").Value = 5
$ws.Range("F458").Value = 100112043
$ws.Range("G458").Value = "Pepino ensalada"
$ws.Range("H458").Value = "Sin especificar"
$ws.Range("I458").Value = "Primera"
$ws.Range("J458").Value = 87
$ws.Range("K458").Value = 17000
$ws.Range("L458").Value = 18000
$ws.Range("M458").Value = 17517
$ws.Range("N458").Value = "$/caja 60 unidades"
$ws.Range("O458").Value = "Limache"
$ws.Range("P458").Value = 292
$ws.Range("Q458").Value = 60
$ws.Range("R458").Value = "Hortaliza"

# D column (Fecha) uses a date numeric format across the table; make sure the
# newly inserted row keeps the same number format as the rest of column D.
$ws.Range("D458").NumberFormat = $ws.Range("D459").NumberFormat
